# edit.ps1 - "Altura de peças padrão automaticamente"
#
# The source system renumbered orders: what used to be order 250209_0001
# (2 pieces) and 250209_0002 (1 piece) at the top of the sheet were dropped,
# every later order shifted up by one order-slot, and new orders
# 250211_0006 through 250212_0005 were appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first three data rows (old orders 250209_0001 / 250209_0002).
# Excel shifts every row below up, so former row 5 becomes row 2, etc.
$ws.Range("A2:A4").EntireRow.Delete()

# Append the new trailing orders (new rows 33-40).

# Row 33: 250211_0006 / 250211_0006_001
$ws.Cells.Item(33, 1).Value = "250211_0006"
$ws.Cells.Item(33, 2).Value = "250211_0006_001"
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 4).Value = "DOUGLAS"
$ws.Cells.Item(33, 5).Value = 7
$ws.Cells.Item(33, 6).Value = 11
$ws.Cells.Item(33, 7).Value = "Peça Fixa"
$ws.Cells.Item(33, 8).Value = 1
$ws.Cells.Item(33, 9).Value = 1000
$ws.Cells.Item(33, 10).Value = 1000
$ws.Cells.Item(33, 11).Value = 975
$ws.Cells.Item(33, 12).Value = 500
$ws.Cells.Item(33, 13).Value = 0.5
$ws.Cells.Item(33, 14).Value = 341.82
$ws.Cells.Item(33, 15).Value = 170.91
$ws.Cells.Item(33, 16).NumberFormat = "@"  # keep "123" as text, not a number
$ws.Cells.Item(33, 16).Value = "123"

# Row 34: 250211_0006 / 250211_0006_002
$ws.Cells.Item(34, 1).Value = "250211_0006"
$ws.Cells.Item(34, 2).Value = "250211_0006_002"
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = "DOUGLAS"
$ws.Cells.Item(34, 5).Value = 7
$ws.Cells.Item(34, 6).Value = 11
$ws.Cells.Item(34, 7).Value = "Peça Móvel"
$ws.Cells.Item(34, 8).Value = 1
$ws.Cells.Item(34, 9).Value = 1000
$ws.Cells.Item(34, 10).Value = 1000
$ws.Cells.Item(34, 11).Value = 938
$ws.Cells.Item(34, 12).Value = 550
$ws.Cells.Item(34, 13).Value = 0.75
$ws.Cells.Item(34, 14).Value = 341.82
$ws.Cells.Item(34, 15).Value = 256.36
$ws.Cells.Item(34, 16).NumberFormat = "@"  # keep "123" as text, not a number
$ws.Cells.Item(34, 16).Value = "123"

# Row 35: 250212_0001 / 250212_0001_001
$ws.Cells.Item(35, 1).Value = "250212_0001"
$ws.Cells.Item(35, 2).Value = "250212_0001_001"
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(35, 4).Value = "DOUGLAS"
$ws.Cells.Item(35, 5).Value = 7
$ws.Cells.Item(35, 6).Value = 11
$ws.Cells.Item(35, 7).Value = "Peça Fixa"
$ws.Cells.Item(35, 8).Value = 1
$ws.Cells.Item(35, 9).Value = 1000
$ws.Cells.Item(35, 10).Value = 1000
$ws.Cells.Item(35, 11).Value = 975
$ws.Cells.Item(35, 12).Value = 500
$ws.Cells.Item(35, 13).Value = 0.5
$ws.Cells.Item(35, 14).Value = 341.82
$ws.Cells.Item(35, 15).Value = 170.91
$ws.Cells.Item(35, 16).NumberFormat = "@"  # keep "47845" as text, not a number
$ws.Cells.Item(35, 16).Value = "47845"

# Row 36: 250212_0001 / 250212_0001_002
$ws.Cells.Item(36, 1).Value = "250212_0001"
$ws.Cells.Item(36, 2).Value = "250212_0001_002"
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(36, 4).Value = "DOUGLAS"
$ws.Cells.Item(36, 5).Value = 7
$ws.Cells.Item(36, 6).Value = 11
$ws.Cells.Item(36, 7).Value = "Peça Móvel"
$ws.Cells.Item(36, 8).Value = 1
$ws.Cells.Item(36, 9).Value = 1000
$ws.Cells.Item(36, 10).Value = 1000
$ws.Cells.Item(36, 11).Value = 938
$ws.Cells.Item(36, 12).Value = 550
$ws.Cells.Item(36, 13).Value = 0.75
$ws.Cells.Item(36, 14).Value = 341.82
$ws.Cells.Item(36, 15).Value = 256.36
$ws.Cells.Item(36, 16).NumberFormat = "@"  # keep "47845" as text, not a number
$ws.Cells.Item(36, 16).Value = "47845"

# Row 37: 250212_0002 / 250212_0002_001
$ws.Cells.Item(37, 1).Value = "250212_0002"
$ws.Cells.Item(37, 2).Value = "250212_0002_001"
$ws.Cells.Item(37, 3).Value = 1
$ws.Cells.Item(37, 4).Value = "DOUGLAS"
$ws.Cells.Item(37, 5).Value = 28
$ws.Cells.Item(37, 6).Value = 2
$ws.Cells.Item(37, 7).Value = "Peça Principal"
$ws.Cells.Item(37, 8).Value = 1
$ws.Cells.Item(37, 9).Value = 1845
$ws.Cells.Item(37, 10).Value = 700
$ws.Cells.Item(37, 11).Value = 1845
$ws.Cells.Item(37, 12).Value = 700
$ws.Cells.Item(37, 13).Value = 1.5
$ws.Cells.Item(37, 14).Value = 205.75
$ws.Cells.Item(37, 15).Value = 308.62
$ws.Cells.Item(37, 16).NumberFormat = "@"  # keep "47845" as text, not a number
$ws.Cells.Item(37, 16).Value = "47845"

# Row 38: 250212_0003 / 250212_0003_001
$ws.Cells.Item(38, 1).Value = "250212_0003"
$ws.Cells.Item(38, 2).Value = "250212_0003_001"
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(38, 4).Value = "DOUGLAS"
$ws.Cells.Item(38, 5).Value = 68
$ws.Cells.Item(38, 6).Value = 20
$ws.Cells.Item(38, 7).Value = "Peça Principal"
$ws.Cells.Item(38, 8).Value = 3
$ws.Cells.Item(38, 9).Value = 975
$ws.Cells.Item(38, 10).Value = 600
$ws.Cells.Item(38, 11).Value = 975
$ws.Cells.Item(38, 12).Value = 600
$ws.Cells.Item(38, 13).Value = 2
$ws.Cells.Item(38, 14).Value = 332.75
$ws.Cells.Item(38, 15).Value = 665.5
$ws.Cells.Item(38, 16).NumberFormat = "@"  # keep "47845" as text, not a number
$ws.Cells.Item(38, 16).Value = "47845"

# Row 39: 250212_0004 / 250212_0004_001
$ws.Cells.Item(39, 1).Value = "250212_0004"
$ws.Cells.Item(39, 2).Value = "250212_0004_001"
$ws.Cells.Item(39, 3).Value = 1
$ws.Cells.Item(39, 4).Value = "DOUGLAS"
$ws.Cells.Item(39, 5).Value = 1
$ws.Cells.Item(39, 6).Value = 2
$ws.Cells.Item(39, 7).Value = "Peça Principal"
$ws.Cells.Item(39, 8).Value = 5
$ws.Cells.Item(39, 9).Value = 980
$ws.Cells.Item(39, 10).Value = 475
$ws.Cells.Item(39, 11).Value = 980
$ws.Cells.Item(39, 12).Value = 475
$ws.Cells.Item(39, 13).Value = 2.5
$ws.Cells.Item(39, 14).Value = 205.75
$ws.Cells.Item(39, 15).Value = 514.38
$ws.Cells.Item(39, 16).NumberFormat = "@"  # keep "47845" as text, not a number
$ws.Cells.Item(39, 16).Value = "47845"

# Row 40: 250212_0005 / 250212_0005_001
$ws.Cells.Item(40, 1).Value = "250212_0005"
$ws.Cells.Item(40, 2).Value = "250212_0005_001"
$ws.Cells.Item(40, 3).Value = 1
$ws.Cells.Item(40, 4).Value = "DOUGLAS"
$ws.Cells.Item(40, 5).Value = 12
$ws.Cells.Item(40, 6).Value = 12
$ws.Cells.Item(40, 7).Value = "Peça Principal"
$ws.Cells.Item(40, 8).Value = 3
$ws.Cells.Item(40, 9).Value = 500
$ws.Cells.Item(40, 10).Value = 300
$ws.Cells.Item(40, 11).Value = 500
$ws.Cells.Item(40, 12).Value = 300
$ws.Cells.Item(40, 13).Value = 0.5
$ws.Cells.Item(40, 14).Value = 483.79
$ws.Cells.Item(40, 15).Value = 241.9
$ws.Cells.Item(40, 16).Value = "Faue"
